# Regenerate the "K" column (column G) values on Sheet1.
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" - the visible effect in this workbook is that the
# per-row K values (column G, header "K") are recalculated and rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2 through 62 (column G), in row order.
$kValues = @(
    0, 1, 2, 0, 0, 1, 1, 0, 0, 0,
    2, 1, 1, 0, 1, 3, 0, 1, 0, 0,
    0, 0, 1, 3, 1, 1, 2, 0, 0, 1,
    0, 1, 0, 1, 1, 0, 1, 0, 1, 1,
    1, 1, 2, 3, 1, 1, 0, 4, 2, 0,
    0, 0, 0, 1, 0, 2, 0, 3, 3, 0,
    2
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
